$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 updates ---
$ws.Range("L6").Value = "18"
$ws.Range("V6").Value = "ui/assets/topdown/top-down-shooter/characters/body/3.png"

# --- Row 7 updates ---
$ws.Range("V7").Value = "ui/assets/topdown/top-down-shooter/characters/body/2.png"
$ws.Range("W7").Value = "0.94"

# --- Row 8 updates ---
$ws.Range("H8").Value = "3.6"
$ws.Range("J8").Value = "BURST"
$ws.Range("K8").Value = "2.60"
$ws.Range("L8").Value = "28"
$ws.Range("M8").Value = "0.70"
$ws.Range("V8").Value = "ui/assets/topdown/top-down-shooter/characters/body/1.png"
$ws.Range("W8").Value = "0.92"

# --- Row 9 updates ---
$ws.Range("L9").Value = "64"
$ws.Range("M9").Value = "0.60"
$ws.Range("V9").Value = "ui/assets/topdown/top-down-shooter/characters/turret/2.png"

# --- Row 11 updates ---
$ws.Range("F11").Value = "130"
$ws.Range("G11").Value = "14"
$ws.Range("H11").Value = "3.6"
$ws.Range("K11").Value = "2.40"
$ws.Range("L11").Value = "18"
$ws.Range("M11").Value = "0.90"
$ws.Range("S11").Value = "6"
$ws.Range("T11").Value = "碎影快速游走射击，首波练习走位。"
$ws.Range("U11").Value = "16"
$ws.Range("V11").Value = "ui/assets/topdown/top-down-shooter/characters/head/9.png"
$ws.Range("W11").Value = "0.8"

# --- New row 12 ---
$ws.Range("A12").Value = "40"
$ws.Range("B12").Value = "06"
$ws.Range("C12").Value = "0007"
$ws.Range("D12").Value = "虚潮行者"
$ws.Range("E12").Value = "SHAMBLER"
$ws.Range("F12").Value = "150"
$ws.Range("G12").Value = "12"
$ws.Range("H12").Value = "2.9"
$ws.Range("I12").Value = "16"
$ws.Range("J12").Value = "MANUAL"
$ws.Range("K12").Value = "0"
$ws.Range("L12").Value = "0"
$ws.Range("M12").Value = "0"
$ws.Range("N12").Value = ""
$ws.Range("O12").Value = ""
$ws.Range("P12").Value = "LIGHT"
$ws.Range("Q12").Value = "VOID"
$ws.Range("R12").Value = "loot:ichor_minor"
$ws.Range("S12").Value = "4"
$ws.Range("T12").Value = "只会贴身缠斗的虚潮行者，用来熟悉位移。"
$ws.Range("U12").Value = "14"
$ws.Range("V12").Value = "ui/assets/topdown/top-down-shooter/characters/head/5.png"
$ws.Range("W12").Value = "0.88"
$ws.Range("X12").Value = "ui/assets/topdown/top-down-shooter/effects/explosion.png"
$ws.Range("Y12").Value = "ui/assets/topdown/top-down-shooter/sounds/death.wav"
$ws.Range("Z12").Value = ""
